# Update column F ("dSF") values for specific rows per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = -5
    6  = -4
    11 = -9
    14 = -4
    16 = 0
    25 = 0
    26 = -3
    28 = 0
    30 = 0
    31 = 4
    32 = -4
    33 = 1
    38 = -4
    39 = -1
    40 = -1
    42 = 5
    44 = 3
    47 = -1
    50 = 2
    52 = 0
    53 = 0
    58 = 8
    59 = 1
    62 = -4
    63 = 3
    64 = -8
    68 = 0
    69 = -1
    70 = -5
    77 = 5
    81 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
